$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R: year 2021 header, matching the style of Q4 (the previous last year)
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# New data point for 2021, matching the style of Q5 but with a dedicated "0.0" number format
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Move the active selection as recorded in the sheet view
$ws.Range("S9").Select()
